$wb = $excel.ActiveWorkbook

# --- Rename existing "3_14" sheet to "3_14_N2" ---
$s3 = $wb.Worksheets.Item("3_14")
$s3.Name = "3_14_N2"

# --- Add the three new worksheets, in order, right after 3_14_N2 ---
$s4 = $wb.Worksheets.Add($null, $s3)
$s4.Name = "3_14_Ar"

$s5 = $wb.Worksheets.Add($null, $s4)
$s5.Name = "1_223_kerogen_N2"

$s6 = $wb.Worksheets.Add($null, $s5)
$s6.Name = "1_223_kerogen_Ar"

# --- Headers (shared strings 22 = "D", 23 = "IPV") for each new sheet ---
$s4.Cells.Item(1, 1).Value = "D"
$s4.Cells.Item(1, 2).Value = "IPV"
$s5.Cells.Item(1, 1).Value = "D"
$s5.Cells.Item(1, 2).Value = "IPV"
$s6.Cells.Item(1, 1).Value = "D"
$s6.Cells.Item(1, 2).Value = "IPV"

# --- Data rows ---

# Sheet 3_14_Ar data
$s4.Cells.Item(2, 1).Value = 2526.1127212088631
$s4.Cells.Item(2, 2).Value = 0.02248987212767856
$s4.Cells.Item(3, 1).Value = 1736.684961001795
$s4.Cells.Item(3, 2).Value = 0.024369528089523588
$s4.Cells.Item(4, 1).Value = 1113.664757437163
$s4.Cells.Item(4, 2).Value = 0.028217193600854767
$s4.Cells.Item(5, 1).Value = 776.84788690406492
$s4.Cells.Item(5, 2).Value = 0.018655220382101712
$s4.Cells.Item(6, 1).Value = 577.74235763144191
$s4.Cells.Item(6, 2).Value = 0.014880110188305715
$s4.Cells.Item(7, 1).Value = 451.37501745590112
$s4.Cells.Item(7, 2).Value = 0.01137179802270074
$s4.Cells.Item(8, 1).Value = 328.43084728994631
$s4.Cells.Item(8, 2).Value = 0.015556032565181004
$s4.Cells.Item(9, 1).Value = 246.76871562327369
$s4.Cells.Item(9, 2).Value = 0.010810974079368128
$s4.Cells.Item(10, 1).Value = 196.42785986524063
$s4.Cells.Item(10, 2).Value = 0.0083128284938789303
$s4.Cells.Item(11, 1).Value = 163.29400530057268
$s4.Cells.Item(11, 2).Value = 0.0063505635124616556
$s4.Cells.Item(12, 1).Value = 138.99052711583241
$s4.Cells.Item(12, 2).Value = 0.0054950298405270666
$s4.Cells.Item(13, 1).Value = 120.78113962659299
$s4.Cells.Item(13, 2).Value = 0.0045992471303288917
$s4.Cells.Item(14, 1).Value = 99.836404890283575
$s4.Cells.Item(14, 2).Value = 0.0067076829289783946
$s4.Cells.Item(15, 1).Value = 81.835127138862418
$s4.Cells.Item(15, 2).Value = 0.0050701804347551924
$s4.Cells.Item(16, 1).Value = 69.046853031696685
$s4.Cells.Item(16, 2).Value = 0.0041692908299011382
$s4.Cells.Item(17, 1).Value = 59.381538225515435
$s4.Cells.Item(17, 2).Value = 0.0035422541980601741
$s4.Cells.Item(18, 1).Value = 51.740425820280016
$s4.Cells.Item(18, 2).Value = 0.0030014598003719642
$s4.Cells.Item(19, 1).Value = 45.516978062998383
$s4.Cells.Item(19, 2).Value = 0.0026759384224303616
$s4.Cells.Item(20, 1).Value = 40.352748402259714
$s4.Cells.Item(20, 2).Value = 0.0024122568929177894
$s4.Cells.Item(21, 1).Value = 35.996420638346237
$s4.Cells.Item(21, 2).Value = 0.0023968028744053606
$s4.Cells.Item(22, 1).Value = 32.228951499691561
$s4.Cells.Item(22, 2).Value = 0.0021626510566854443
$s4.Cells.Item(23, 1).Value = 28.81486187018913
$s4.Cells.Item(23, 2).Value = 0.0022844694920196324
$s4.Cells.Item(24, 1).Value = 25.730783184562011
$s4.Cells.Item(24, 2).Value = 0.0022205885899849024
$s4.Cells.Item(25, 1).Value = 22.893579080704612
$s4.Cells.Item(25, 2).Value = 0.002147992385151112
$s4.Cells.Item(26, 1).Value = 20.212302104662434
$s4.Cells.Item(26, 2).Value = 0.0025232743643620998

# Sheet 1_223_kerogen_N2 data
$s5.Cells.Item(2, 1).Value = 2522.6986019665942
$s5.Cells.Item(2, 2).Value = 0.012982483065113449
$s5.Cells.Item(3, 1).Value = 1819.3653643428077
$s5.Cells.Item(3, 2).Value = 0.084082202822941057
$s5.Cells.Item(4, 1).Value = 1093.7207752333341
$s5.Cells.Item(4, 2).Value = 0.084562256166683364
$s5.Cells.Item(5, 1).Value = 816.8448309357251
$s5.Cells.Item(5, 2).Value = 0.030244054381522407
$s5.Cells.Item(6, 1).Value = 238.35465918336786
$s5.Cells.Item(6, 2).Value = 0.15655967737395926
$s5.Cells.Item(7, 1).Value = 196.46497313182095
$s5.Cells.Item(7, 2).Value = 0.0067456155534584708
$s5.Cells.Item(8, 1).Value = 182.22501431672362
$s5.Cells.Item(8, 2).Value = 0.0059840523354105085
$s5.Cells.Item(9, 1).Value = 169.33633364128332
$s5.Cells.Item(9, 2).Value = 0.0053524022497884755
$s5.Cells.Item(10, 1).Value = 158.15074059204255
$s5.Cells.Item(10, 2).Value = 0.0053768722724629086
$s5.Cells.Item(11, 1).Value = 148.36049511470196
$s5.Cells.Item(11, 2).Value = 0.0049181297752207073
$s5.Cells.Item(12, 1).Value = 139.25174603774286
$s5.Cells.Item(12, 2).Value = 0.0042998303256459999
$s5.Cells.Item(13, 1).Value = 130.39320225088008
$s5.Cells.Item(13, 2).Value = 0.0043574265827719359
$s5.Cells.Item(14, 1).Value = 122.50261452783785
$s5.Cells.Item(14, 2).Value = 0.0038109283093383829
$s5.Cells.Item(15, 1).Value = 115.76441637384316
$s5.Cells.Item(15, 2).Value = 0.0039495178600387178
$s5.Cells.Item(16, 1).Value = 109.55153675205362
$s5.Cells.Item(16, 2).Value = 0.0034389737096424226
$s5.Cells.Item(17, 1).Value = 103.85320460819905
$s5.Cells.Item(17, 2).Value = 0.0035517506008682139
$s5.Cells.Item(18, 1).Value = 98.531459223327531
$s5.Cells.Item(18, 2).Value = 0.0028762871945155594
$s5.Cells.Item(19, 1).Value = 93.56770891810315
$s5.Cells.Item(19, 2).Value = 0.0030440630420086748
$s5.Cells.Item(20, 1).Value = 89.068978048778504
$s5.Cells.Item(20, 2).Value = 0.0028694208719804125
$s5.Cells.Item(21, 1).Value = 84.927713001740685
$s5.Cells.Item(21, 2).Value = 0.0028236033428654457
$s5.Cells.Item(22, 1).Value = 81.041101044566872
$s5.Cells.Item(22, 2).Value = 0.0027255784731393729
$s5.Cells.Item(23, 1).Value = 77.36307713796846
$s5.Cells.Item(23, 2).Value = 0.0021923485838683614
$s5.Cells.Item(24, 1).Value = 73.93827961554679
$s5.Cells.Item(24, 2).Value = 0.002075961024126838
$s5.Cells.Item(25, 1).Value = 70.628498942275982
$s5.Cells.Item(25, 2).Value = 0.0023471484133108554
$s5.Cells.Item(26, 1).Value = 67.420994584158151
$s5.Cells.Item(26, 2).Value = 0.0020486634369024669
$s5.Cells.Item(27, 1).Value = 64.398207018287422
$s5.Cells.Item(27, 2).Value = 0.0019931035791974127
$s5.Cells.Item(28, 1).Value = 61.672616933136204
$s5.Cells.Item(28, 2).Value = 0.0021197769239574451
$s5.Cells.Item(29, 1).Value = 59.115343085250039
$s5.Cells.Item(29, 2).Value = 0.0020160347501587796
$s5.Cells.Item(30, 1).Value = 56.708636658171834
$s5.Cells.Item(30, 2).Value = 0.0018912198159500889
$s5.Cells.Item(31, 1).Value = 54.462479222632375
$s5.Cells.Item(31, 2).Value = 0.0019085335344718871
$s5.Cells.Item(32, 1).Value = 52.288087554525781
$s5.Cells.Item(32, 2).Value = 0.0017973681202013078
$s5.Cells.Item(33, 1).Value = 50.245756915724122
$s5.Cells.Item(33, 2).Value = 0.0017078139876405825
$s5.Cells.Item(34, 1).Value = 45.569630404997604
$s5.Cells.Item(34, 2).Value = 0.0044245882830829921
$s5.Cells.Item(35, 1).Value = 39.840598529263879
$s5.Cells.Item(35, 2).Value = 0.00408348954030055
$s5.Cells.Item(36, 1).Value = 35.108574573019943
$s5.Cells.Item(36, 2).Value = 0.0033054588816452644
$s5.Cells.Item(37, 1).Value = 31.063274620405899
$s5.Cells.Item(37, 2).Value = 0.0030387253466121213
$s5.Cells.Item(38, 1).Value = 27.596139215799056
$s5.Cells.Item(38, 2).Value = 0.0027200859448235446
$s5.Cells.Item(39, 1).Value = 24.461259512273841
$s5.Cells.Item(39, 2).Value = 0.0029558925936602064
$s5.Cells.Item(40, 1).Value = 21.526187504475047
$s5.Cells.Item(40, 2).Value = 0.0031887114289225792
$s5.Cells.Item(41, 1).Value = 18.884628381438144
$s5.Cells.Item(41, 2).Value = 0.0026478492972980415
$s5.Cells.Item(42, 1).Value = 17.283311554507677
$s5.Cells.Item(42, 2).Value = 0.0014691730063481556

# Sheet 1_223_kerogen_Ar data
$s6.Cells.Item(2, 1).Value = 2546.3780987380183
$s6.Cells.Item(2, 2).Value = 0.023514104997692166
$s6.Cells.Item(3, 1).Value = 2148.8454639800379
$s6.Cells.Item(3, 2).Value = 0.023889537108014877
$s6.Cells.Item(4, 1).Value = 1547.9669862561786
$s6.Cells.Item(4, 2).Value = 0.05311778713870377
$s6.Cells.Item(5, 1).Value = 991.11292174657683
$s6.Cells.Item(5, 2).Value = 0.06018146766059046
$s6.Cells.Item(6, 1).Value = 722.32
$s6.Cells.Item(6, 2).Value = 0.032699219670457823
$s6.Cells.Item(7, 1).Value = 569.5442221574051
$s6.Cells.Item(7, 2).Value = 0.028323442779093332
$s6.Cells.Item(8, 1).Value = 463.78135802977192
$s6.Cells.Item(8, 2).Value = 0.021911648301101209
$s6.Cells.Item(9, 1).Value = 323.96873387653483
$s6.Cells.Item(9, 2).Value = 0.044596808674386472
$s6.Cells.Item(10, 1).Value = 240.74462929968516
$s6.Cells.Item(10, 2).Value = 0.024531055901001523
$s6.Cells.Item(11, 1).Value = 193.19157974600617
$s6.Cells.Item(11, 2).Value = 0.018552544802818552
$s6.Cells.Item(12, 1).Value = 161.43447788306068
$s6.Cells.Item(12, 2).Value = 0.014360692859056581
$s6.Cells.Item(13, 1).Value = 138.06087926736436
$s6.Cells.Item(13, 2).Value = 0.012149941074725187
$s6.Cells.Item(14, 1).Value = 119.66974559385622
$s6.Cells.Item(14, 2).Value = 0.010209658072152057
$s6.Cells.Item(15, 1).Value = 98.893108389781588
$s6.Cells.Item(15, 2).Value = 0.014959830777574228
$s6.Cells.Item(16, 1).Value = 81.290143928126156
$s6.Cells.Item(16, 2).Value = 0.011374296515418734
$s6.Cells.Item(17, 1).Value = 68.681444654955286
$s6.Cells.Item(17, 2).Value = 0.0091219505445371881
$s6.Cells.Item(18, 1).Value = 59.152364069403546
$s6.Cells.Item(18, 2).Value = 0.007436682736656295
$s6.Cells.Item(19, 1).Value = 51.582718869364172
$s6.Cells.Item(19, 2).Value = 0.0068615094496146381
$s6.Cells.Item(20, 1).Value = 45.418234690783805
$s6.Cells.Item(20, 2).Value = 0.0059148325627828378
$s6.Cells.Item(21, 1).Value = 40.319771141873289
$s6.Cells.Item(21, 2).Value = 0.0054687141862738396
$s6.Cells.Item(22, 1).Value = 35.937134694405081
$s6.Cells.Item(22, 2).Value = 0.0052245541073886285
$s6.Cells.Item(23, 1).Value = 32.21175362046661
$s6.Cells.Item(23, 2).Value = 0.0048941050703838031
$s6.Cells.Item(24, 1).Value = 28.822655254359148
$s6.Cells.Item(24, 2).Value = 0.0053232805694557643
$s6.Cells.Item(25, 1).Value = 25.699190706885197
$s6.Cells.Item(25, 2).Value = 0.0056101765541920435
$s6.Cells.Item(26, 1).Value = 22.872788963220255
$s6.Cells.Item(26, 2).Value = 0.0056522394732331496
$s6.Cells.Item(27, 1).Value = 20.158758779355104
$s6.Cells.Item(27, 2).Value = 0.0064831588583638411

# --- Selections: sheets 3_14_N2 / 3_14_Ar / 1_223_kerogen_N2 select A1:B1,
#     with no cell actively edited; 1_223_kerogen_Ar (last-added, now the
#     active tab) keeps the default active-cell selection and becomes the
#     selected/visible tab. ---
$s3.Activate()
[void]$s3.Range("A1:B1").Select()

$s4.Activate()
[void]$s4.Range("A1:B1").Select()

$s5.Activate()
[void]$s5.Range("A1:B1").Select()

$s6.Activate()
[void]$s6.Range("F9").Select()
